# Apply final correction edits to the "list categories" sheet of the
# FSD11k dataset draft workbook (Chime & Finger snapping categories
# sent to dur split; totals / counts updated accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list categories")

# Overall total number of sounds in final dataset (row 1, col C)
$ws.Range("C1").Value = 11073

# Row 4 (second data category): num eval count decreases by 1
$ws.Range("G4").Value = 25

# Row 34: "Chime" category counts
$ws.Range("C34").Value = 66   # num dev LQ
$ws.Range("F34").Value = 115  # num dev final
$ws.Range("G34").Value = 29   # num eval

# Row 38: "Finger snapping" category counts
$ws.Range("C38").Value = 77   # num dev LQ
$ws.Range("D38").Value = 40   # num dev LQ prior
$ws.Range("F38").Value = 117  # num dev final
$ws.Range("G38").Value = 33   # num eval

$wb.Save()
